$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 60; this pushes the existing rows 60-160
# down to 61-161 (preserving all of their data/styles).
$ws.Rows(60).Insert()

# Populate the newly inserted row 60 with the new weekly record.
$ws.Range("A60").Value = 8
$ws.Range("B60").Value = "Terminal La Palmera de La Serena"
$ws.Range("C60").Value = "Coquimbo"
$ws.Range("D60").Value = 44477
$ws.Range("D60").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E60").Value = 4
$ws.Range("F60").Value = 100112012
$ws.Range("G60").Value = "Espinaca"
$ws.Range("H60").Value = "Sin especificar"
$ws.Range("I60").Value = "Primera"
$ws.Range("J60").Value = 3300
$ws.Range("K60").Value = 400
$ws.Range("L60").Value = 500
$ws.Range("M60").Value = 450
$ws.Range("N60").Value = "`$/atado 300 a 500 gramos"
$ws.Range("O60").Value = "Provincia del Elquí"
$ws.Range("P60").Value = 900
$ws.Range("Q60").Value = 0.5
$ws.Range("R60").Value = "Hortaliza"
